# Add a new "V2G_Scenarios" row to the Index Table on the Setting_V2G_in_EU sheet.
# This mirrors the author's edit: a new row is inserted before the existing
# "Energy_Storage_Scenarios" row (row 29), pushing that row and everything
# below it down by one, and the new row is populated with the V2G aspect data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting_V2G_in_EU")

# Insert a new blank row above row 29 (shifts row 29.. down to 30..)
$ws.Rows("29:29").Insert()

# Populate the new row 29 with the V2G_Scenarios aspect entry
$ws.Range("C29").Value = "V2G_Scenarios"
$ws.Range("D29").Value = "V2G share of vehicle fleet scenarios"
$ws.Range("F29").Value = "V2G_Scenarios"
$ws.Range("G29").Value = "All"
$ws.Range("H29").Value = "v"
$ws.Range("I29").Value = "v2g"

# Match styling used by the neighbouring rows (C has style 45, H has style 28)
$ws.Range("C29").Style = $ws.Range("C30").Style
$ws.Range("H29").Style = $ws.Range("H30").Style

$ws.Range("J29").Select()
